$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data between row 106 and row 107 ---
# (A, B, C, D, E, G stay the same; F, H, I, J, L, M, N, P, Q, R, T, U, V swap)
# NOTE: this engine's Range.Value getter is unreliable, use Value2 to read.
$row106 = @{
    F = $ws.Range("F106").Value2
    H = $ws.Range("H106").Value2
    I = $ws.Range("I106").Value2
    J = $ws.Range("J106").Value2
    L = $ws.Range("L106").Value2
    M = $ws.Range("M106").Value2
    N = $ws.Range("N106").Value2
    P = $ws.Range("P106").Value2
    Q = $ws.Range("Q106").Value2
    R = $ws.Range("R106").Value2
    T = $ws.Range("T106").Value2
    U = $ws.Range("U106").Value2
    V = $ws.Range("V106").Value2
}
$row107 = @{
    F = $ws.Range("F107").Value2
    H = $ws.Range("H107").Value2
    I = $ws.Range("I107").Value2
    J = $ws.Range("J107").Value2
    L = $ws.Range("L107").Value2
    M = $ws.Range("M107").Value2
    N = $ws.Range("N107").Value2
    P = $ws.Range("P107").Value2
    Q = $ws.Range("Q107").Value2
    R = $ws.Range("R107").Value2
    T = $ws.Range("T107").Value2
    U = $ws.Range("U107").Value2
    V = $ws.Range("V107").Value2
}

$ws.Range("F106").Value = $row107.F
$ws.Range("H106").Value = $row107.H
$ws.Range("I106").Value = $row107.I
$ws.Range("J106").Value = $row107.J
$ws.Range("L106").Value = $row107.L
$ws.Range("M106").Value = $row107.M
$ws.Range("N106").Value = $row107.N
$ws.Range("P106").Value = $row107.P
$ws.Range("Q106").Value = $row107.Q
$ws.Range("R106").Value = $row107.R
$ws.Range("T106").Value = $row107.T
$ws.Range("U106").Value = $row107.U
$ws.Range("V106").Value = $row107.V

$ws.Range("F107").Value = $row106.F
$ws.Range("H107").Value = $row106.H
$ws.Range("I107").Value = $row106.I
$ws.Range("J107").Value = $row106.J
$ws.Range("L107").Value = $row106.L
$ws.Range("M107").Value = $row106.M
$ws.Range("N107").Value = $row106.N
$ws.Range("P107").Value = $row106.P
$ws.Range("Q107").Value = $row106.Q
$ws.Range("R107").Value = $row106.R
$ws.Range("T107").Value = $row106.T
$ws.Range("U107").Value = $row106.U
$ws.Range("V107").Value = $row106.V

# --- Append a new match row (row 121) ---
$ws.Range("A121").Value2 = 120
$ws.Range("B121").Value2 = "romania"
$ws.Range("C121").Value2 = "liga-2"
$ws.Range("D121").Value2 = "2023-2024"
$ws.Range("E121").Value2 = 45236.58333333334
$ws.Range("F121").Value2 = "Unirea Slobozia"
$ws.Range("G121").Value2 = 0
$ws.Range("H121").Value2 = "Gloria Buzau"
$ws.Range("I121").Value2 = 1
$ws.Range("J121").Value2 = 2.33
$ws.Range("K121").Value2 = "02/11/2023 22:12"
$ws.Range("L121").Value2 = 2.76
$ws.Range("M121").Value2 = "05/11/2023 09:28"
$ws.Range("N121").Value2 = 2.98
$ws.Range("O121").Value2 = "02/11/2023 22:12"
$ws.Range("P121").Value2 = 3.05
$ws.Range("Q121").Value2 = "05/11/2023 02:59"
$ws.Range("R121").Value2 = 2.96
$ws.Range("S121").Value2 = "02/11/2023 22:12"
$ws.Range("T121").Value2 = 2.67
$ws.Range("U121").Value2 = "05/11/2023 09:28"
$ws.Range("V121").Value2 = "https://www.betexplorer.com/football/romania/liga-2/unirea-slobozia-fc-buzau/Mqj1DJXb/"

# Match the styling used by the other data rows (A col bordered/bold/centered, E col date-time format)
# (direct Style object assignment doesn't stick on this engine, so copy/paste-special formats)
$ws.Range("A120").Copy()
$ws.Range("A121").PasteSpecial(-4122)
$ws.Range("E120").Copy()
$ws.Range("E121").PasteSpecial(-4122)
$excel.CutCopyMode = 0
